# Auto-generated: apply cell value updates per diff
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 783.2157
$ws.Range("J17").Value = 775.8163500000001
$ws.Range("L17").Value = 2327.44905
$ws.Range("N17").Value = -2663.44905
$ws.Range("H40").Value = 2352.4285
$ws.Range("I40").Value = 1928.125
$ws.Range("J40").Value = 2613.5386
$ws.Range("K40").Value = 1928.125
$ws.Range("L40").Value = 2613.5386
$ws.Range("M40").Value = -1753.125
$ws.Range("N40").Value = -2963.5386
$ws.Range("H69").Value = 8559.866
$ws.Range("I69").Value = 5183.3335
$ws.Range("J69").Value = 9404
$ws.Range("K69").Value = 15550.0005
$ws.Range("L69").Value = 28212
$ws.Range("M69").Value = -14676.0005
$ws.Range("N69").Value = -29960
$ws.Range("H72").Value = 8559.866
$ws.Range("I72").Value = 5183.3335
$ws.Range("J72").Value = 9404
$ws.Range("K72").Value = 46650.0015
$ws.Range("L72").Value = 84636
$ws.Range("M72").Value = -42282.0015
$ws.Range("N72").Value = -93372
$ws.Range("H116").Value = 224718.89
$ws.Range("I116").Value = 90177
$ws.Range("J116").Value = 303201.66
$ws.Range("K116").Value = 90177
$ws.Range("L116").Value = 303201.66
$ws.Range("M116").Value = -86735
$ws.Range("N116").Value = -310085.66
$ws.Range("H132").Value = 61771.95
$ws.Range("I132").Value = 68137.375
$ws.Range("K132").Value = 204412.125
$ws.Range("M132").Value = -201882.125
$ws.Range("H135").Value = 1837.826
$ws.Range("I135").Value = 587
$ws.Range("J135").Value = 4183.125
$ws.Range("K135").Value = 5283
$ws.Range("L135").Value = 37648.125
$ws.Range("M135").Value = -2748
$ws.Range("N135").Value = -42718.125
$ws.Range("H137").Value = 339946.75
$ws.Range("I137").Value = 2202.2334
$ws.Range("J137").Value = 493467
$ws.Range("K137").Value = 6606.7002
$ws.Range("L137").Value = 1480401
$ws.Range("M137").Value = -4056.7002
$ws.Range("N137").Value = -1485501
$ws.Range("H138").Value = 1874.5
$ws.Range("I138").Value = 1220.8695
$ws.Range("K138").Value = 3662.6085
$ws.Range("M138").Value = 1477.3915
$ws.Range("H141").Value = 1476.5883
$ws.Range("I141").Value = 1287.6428
$ws.Range("J141").Value = 2358.3333
$ws.Range("K141").Value = 3862.9284
$ws.Range("L141").Value = 7074.999899999999
$ws.Range("M141").Value = 1317.0716
$ws.Range("N141").Value = -17434.9999

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H5").Value = 217
$ws.Range("I5").Value = 146.25
$ws.Range("K5").Value = 146.25
$ws.Range("M5").Value = -34.25
$ws.Range("H32").Value = 12351492
$ws.Range("I32").Value = 12351492
$ws.Range("K32").Value = 12351492
$ws.Range("M32").Value = -12351205
$ws.Range("H61").Value = 4777290
$ws.Range("I61").Value = 4777290
$ws.Range("K61").Value = 4777290
$ws.Range("M61").Value = -4777078
$ws.Range("H74").Value = 10673.353
$ws.Range("I74").Value = 12404.272
$ws.Range("K74").Value = 12404.272
$ws.Range("M74").Value = -11530.272
$ws.Range("H77").Value = 10673.353
$ws.Range("I77").Value = 12404.272
$ws.Range("K77").Value = 62021.36
$ws.Range("M77").Value = -57653.36
$ws.Range("H97").Value = 1786.8889
$ws.Range("I97").Value = 2002.5
$ws.Range("J97").Value = 1614.4
$ws.Range("K97").Value = 2002.5
$ws.Range("L97").Value = 1614.4
$ws.Range("M97").Value = -1506.5
$ws.Range("N97").Value = -2606.4
$ws.Range("H122").Value = 2571.1667
$ws.Range("I122").Value = 2642.625
$ws.Range("K122").Value = 7927.875
$ws.Range("M122").Value = -5477.875
$ws.Range("H125").Value = 0
$ws.Range("I125").Value = 0
$ws.Range("J125").Value = 0
$ws.Range("K125").Value = 0
$ws.Range("H136").Value = 4777290
$ws.Range("I136").Value = 4777290
$ws.Range("K136").Value = 14331870
$ws.Range("M136").Value = -14329320
$ws.Range("L125").ClearContents()
$ws.Range("M125").ClearContents()
$ws.Range("N125").ClearContents()

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 217
$ws.Range("I4").Value = 146.25
$ws.Range("K4").Value = 146.25
$ws.Range("M4").Value = -31.25
$ws.Range("H48").Value = 500342
$ws.Range("J48").Value = 500342
$ws.Range("L48").Value = 500342
$ws.Range("N48").Value = -501172
$ws.Range("H94").Value = 1683.4
$ws.Range("I94").Value = 1351.75
$ws.Range("K94").Value = 1351.75
$ws.Range("M94").Value = -900.75
$ws.Range("H105").Value = 0
$ws.Range("I105").Value = 0
$ws.Range("K105").Value = 0
$ws.Range("H134").Value = 2584553.5
$ws.Range("I134").Value = 3972355.8
$ws.Range("K134").Value = 11917067.4
$ws.Range("M134").Value = -11914532.4
$ws.Range("M105").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H10").Value = 565.2
$ws.Range("I10").Value = 907
$ws.Range("J10").Value = 479.75
$ws.Range("K10").Value = 907
$ws.Range("L10").Value = 479.75
$ws.Range("M10").Value = -768
$ws.Range("N10").Value = -757.75
$ws.Range("H22").Value = 752.06665
$ws.Range("I22").Value = 770.8570999999999
$ws.Range("J22").Value = 489
$ws.Range("K22").Value = 770.8570999999999
$ws.Range("L22").Value = 489
$ws.Range("M22").Value = -420.8570999999999
$ws.Range("N22").Value = -1189
$ws.Range("H41").Value = 23318
$ws.Range("J41").Value = 31421.6
$ws.Range("L41").Value = 31421.6
$ws.Range("N41").Value = -32277.6
$ws.Range("H99").Value = 4322.6
$ws.Range("I99").Value = 3499.5
$ws.Range("K99").Value = 3499.5
$ws.Range("M99").Value = -2001.5
$ws.Range("H126").Value = 4322.6
$ws.Range("I126").Value = 3499.5
$ws.Range("K126").Value = 10498.5
$ws.Range("M126").Value = -8028.5
$ws.Range("H134").Value = 7786.727
$ws.Range("I134").Value = 8621.421
$ws.Range("K134").Value = 25864.263
$ws.Range("M134").Value = -23329.263

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H138").Value = 2983
$ws.Range("I138").Value = 2725
$ws.Range("K138").Value = 8175
$ws.Range("M138").Value = -3035

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 646.4286
$ws.Range("I2").Value = 847.2308
$ws.Range("K2").Value = 847.2308
$ws.Range("M2").Value = -734.2308
$ws.Range("H3").Value = 297.5
$ws.Range("I3").Value = 297.5
$ws.Range("J3").Value = 0
$ws.Range("K3").Value = 297.5
$ws.Range("L3").Value = 0
$ws.Range("H96").Value = 0
$ws.Range("J96").Value = 0
$ws.Range("H122").Value = 55753.25
$ws.Range("I122").Value = 93240.27
$ws.Range("J122").Value = 9935.777
$ws.Range("K122").Value = 279720.81
$ws.Range("L122").Value = 29807.331
$ws.Range("M122").Value = -277270.81
$ws.Range("N122").Value = -34707.331
$ws.Range("M3").ClearContents()
$ws.Range("N3").ClearContents()
$ws.Range("L96").ClearContents()
$ws.Range("N96").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 2349.6667
$ws.Range("I16").Value = 2366
$ws.Range("J16").Value = 2333.3333
$ws.Range("K16").Value = 2366
$ws.Range("L16").Value = 2333.3333
$ws.Range("M16").Value = -2196
$ws.Range("N16").Value = -2673.3333
$ws.Range("H100").Value = 13461.333
$ws.Range("I100").Value = 2643.5
$ws.Range("K100").Value = 2643.5
$ws.Range("M100").Value = -2102.5
$ws.Range("H132").Value = 1202663.2
$ws.Range("I132").Value = 1658459.2
$ws.Range("K132").Value = 4975377.6
$ws.Range("M132").Value = -4972847.6
$ws.Range("H136").Value = 49309
$ws.Range("I136").Value = 2776.375
$ws.Range("K136").Value = 8329.125
$ws.Range("M136").Value = -5779.125

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H46").Value = 89999
$ws.Range("J46").Value = 89999
$ws.Range("L46").Value = 89999
$ws.Range("N46").Value = -90461
$ws.Range("H111").Value = 64322.5
$ws.Range("J111").Value = 64322.5
$ws.Range("L111").Value = 64322.5
$ws.Range("N111").Value = -72502.5
$ws.Range("H113").Value = 3130.9443
$ws.Range("I113").Value = 1276.1428
$ws.Range("J113").Value = 4311.273
$ws.Range("K113").Value = 3828.4284
$ws.Range("L113").Value = 12933.819
$ws.Range("M113").Value = -1658.4284
$ws.Range("N113").Value = -17273.819
$ws.Range("H132").Value = 3663416
$ws.Range("I132").Value = 5298582.5
$ws.Range("J132").Value = 8338.058999999999
$ws.Range("K132").Value = 15895747.5
$ws.Range("L132").Value = 25014.177
$ws.Range("M132").Value = -15893217.5
$ws.Range("N132").Value = -30074.177
$ws.Range("H134").Value = 89999
$ws.Range("J134").Value = 89999
$ws.Range("L134").Value = 269997
$ws.Range("N134").Value = -275067
$ws.Range("H136").Value = 9813.68
$ws.Range("I136").Value = 9697.409
$ws.Range("J136").Value = 10666.333
$ws.Range("K136").Value = 29092.227
$ws.Range("L136").Value = 31998.999
$ws.Range("M136").Value = -26542.227
$ws.Range("N136").Value = -37098.999
